$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the existing layout (title row, old header/columns, stray formatted
# cell at I11) so we can rebuild the compact 4-column table.
$ws.Range("A1:I11").Clear()

# Header row (bold)
$ws.Cells.Item(1,1).Value = "Klinik"
$ws.Cells.Item(1,2).Value = "Adresse"
$ws.Cells.Item(1,3).Value = "PLZ"
$ws.Cells.Item(1,4).Value = "Ort"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true

# Data rows
$ws.Cells.Item(2,1).Value = "Centre hospitalier de Luxembourg"
$ws.Cells.Item(2,2).Value = "Rue Nicolas Ernest Barblé 4"
$ws.Cells.Item(2,3).Value = 1210
$ws.Cells.Item(2,4).Value = "Luxembourg"

$ws.Cells.Item(3,1).Value = "Centre hospitaliser Emile Mayrisch"
$ws.Cells.Item(3,2).Value = "Rue Emile Mayrisch"
$ws.Cells.Item(3,3).Value = 4240
$ws.Cells.Item(3,4).Value = "Esch-sur-Alzette"

$ws.Cells.Item(4,1).Value = "Centre hospitalier du Nord"
$ws.Cells.Item(4,2).Value = "Av. Lucien Salentiny 120"
$ws.Cells.Item(4,3).Value = 9080
$ws.Cells.Item(4,4).Value = "Ettelbruck"

$ws.Cells.Item(5,1).Value = "Hôpital Kirchberg "
$ws.Cells.Item(5,2).Value = "Rue Edward Steichen 9"
$ws.Cells.Item(5,3).Value = 2540
$ws.Cells.Item(5,4).Value = "Neudorf-Weimershof"

# Column widths roughly matching the new, narrower layout.
$ws.Columns.Item(1).ColumnWidth = 31.17
$ws.Columns.Item(2).ColumnWidth = 22.83
$ws.Columns.Item(4).ColumnWidth = 13.83

[void]$ws.Range("E8").Select()
